$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

# Widen column F (Function) from 11 to 20.
# Note: Excel's ColumnWidth (character units) round-trips through a
# pixel-grid conversion when persisted back to the OOXML `width` attribute
# (width = ColumnWidth + 5/6 here). Compensate so the saved width is
# exactly 20.
$ws.Columns.Item(6).ColumnWidth = (20 - 5/6)

# Update the Function value for REQ_003 (row 4) from Clear_Table to Add_Included_Dessert
$ws.Range("F4").Value = "Add_Included_Dessert"
